$wb = $excel.ActiveWorkbook

# --- CapEx sheet: update renewal cost figures ---
$wsCapEx = $wb.Worksheets.Item("CapEx")
$wsCapEx.Range("B2").Value = 500000
$wsCapEx.Range("B3").Value = 1000000
$wsCapEx.Range("B4").Formula = "=8*Identification!C11"

# --- PermanentLoan sheet: update LTV, IO period, amort, and fees ---
$wsPerm = $wb.Worksheets.Item("PermanentLoan")
$wsPerm.Range("C3").Value = 0.65
$wsPerm.Range("D4").Value = 30
$wsPerm.Range("D5").ClearContents()
$wsPerm.Range("C6").Value = 0.01

# --- Update view/selection state to match final UI interaction ---
# User finishes on CapEx at B7, then switches to PermanentLoan landing on A6
$wsCapEx.Activate()
$wsCapEx.Range("B7").Select()

$wsPerm.Activate()
$wsPerm.Range("A6").Select()
